$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in K18 (was missing) on the existing row 18 (run_id 17) ---
$ws.Range("K18").Value = "17%`nopp_won=5%"
$ws.Range("K18").NumberFormat = "0%"

# --- New data rows 19-26 (run_id 18-25) ---

# Row 19 (run_id 18)
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "ppo"
$ws.Range("C19").Value = 128
$ws.Range("D19").Value = 50000
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 5000000
$ws.Range("H19").Value = "change position"
$ws.Range("I19").Value = "ball position`nracket position"
$ws.Range("J19").Value = "hit = 10`nhit opponent court = 10`nelse = 0"
$ws.Range("K19").Value = "27%`nopp_won=5%"
$ws.Range("K19").NumberFormat = "0%"
$ws.Range("L19").Value = "equalizing the reward function for hit and court"

# Row 20 (run_id 19)
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "ppo"
$ws.Range("C20").Value = 128
$ws.Range("D20").Value = 50000
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 5000000
$ws.Range("H20").Value = "rotate angles`n& change position"
$ws.Range("I20").Value = "ball position`nracket position"
$ws.Range("J20").Value = "hit = 10`nhit opponent court = 10`nelse = 0"
$ws.Range("K20").Value = "36%`nopp_won=5%"
$ws.Range("K20").NumberFormat = "0%"
$ws.Range("L20").Value = "add racket rotation to observation and actions"

# Row 21 (run_id 20)
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "ppo"
$ws.Range("C21").Value = 128
$ws.Range("D21").Value = 50000
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 8
$ws.Range("G21").Value = 5000000
$ws.Range("H21").Value = "rotate angles`n& change position"
$ws.Range("I21").Value = "ball position`nracket position"
$ws.Range("J21").Value = "hit = 10`nhit opponent court = 10`nelse = 0"
$ws.Range("K21").Value = "46%`nopp_won=5%"
$ws.Range("K21").NumberFormat = "0%"
$ws.Range("L21").Value = "change stack vector to 15"

# Row 22 (run_id 21)
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "ppo"
$ws.Range("C22").Value = 128
$ws.Range("D22").Value = 50000
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 128
$ws.Range("G22").Value = 5000000
$ws.Range("H22").Value = "rotate angles`n& change position"
$ws.Range("I22").Value = "ball position`nracket position"
$ws.Range("J22").Value = "hit = 10`nhit opponent court = 10`nelse = 0"
$ws.Range("K22").Value = "56%`nopp_won=10%"
$ws.Range("K22").NumberFormat = "0%"
$ws.Range("L22").Value = "bigger network"

# Row 23 (run_id 22)
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "ppo"
$ws.Range("C23").Value = 128
$ws.Range("D23").Value = 50000
$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 128
$ws.Range("G23").Value = 5000000
$ws.Range("H23").Value = "rotate angles`n& change position"
$ws.Range("I23").Value = "ball position`nracket position"
$ws.Range("J23").Value = "hit = 10`nhit opponent court = 10`nelse = 0"
$ws.Range("K23").Value = "45%`nopp_won=7%"
$ws.Range("K23").NumberFormat = "0%"
$ws.Range("L23").Value = "Solved first hit problem"

# Row 24 (run_id 23)
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "ppo"
$ws.Range("C24").Value = 128
$ws.Range("D24").Value = 50000
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 128
$ws.Range("G24").Value = 5000000
$ws.Range("H24").Value = "rotate angles`n& change position"
$ws.Range("I24").Value = "ball position`nracket position"
$ws.Range("J24").Value = "hit = 10`nhit opponent court = 10`nelse = 0"
$ws.Range("K24").Value = "?%`nopp_won=?%"
$ws.Range("K24").NumberFormat = "0%"
$ws.Range("L24").Value = "adding rotation to observation"

# Row 25 (run_id 24)
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "ppo"
$ws.Range("C25").Value = 128
$ws.Range("D25").Value = 50000
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 128
$ws.Range("G25").Value = 5000000
$ws.Range("H25").Value = "rotate angles`n& change position"
$ws.Range("I25").Value = "ball position`nracket position"
$ws.Range("J25").Value = "hit = 1`nhit opponent court = 10`nelse =01"
$ws.Range("K25").Value = "5%`nopp_won=5%"
$ws.Range("K25").NumberFormat = "0%"
$ws.Range("L25").Value = "change reward"

# Row 26 (run_id 25) -- no K value for this row
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "ppo"
$ws.Range("C26").Value = 128
$ws.Range("D26").Value = 50000
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 128
$ws.Range("G26").Value = 5000000
$ws.Range("H26").Value = "rotate angles`n& change position"
$ws.Range("I26").Value = "ball position`nracket position"
$ws.Range("J26").Value = "hit = 10`nhit opponent court = 10`nelse = 0"
$ws.Range("L26").Value = "add speed to observation & change decision period to 4"

# Row 27 (run_id 26) -- only L has a value
$ws.Range("L27").Value = "dynamic reward for hit and court"

# The populated rows (18 through 26) wrap onto multiple lines, so they grow
# to a 60pt row height in the authored file; row 27 keeps the default height.
$ws.Rows("19:26").RowHeight = 60

# Row 1 header reverts to the default (no explicit row height) after resave
$ws.Rows(1).AutoFit() | Out-Null

# Update the active cell/selection to match the authored file's cursor position
$ws.Range("L26").Select() | Out-Null

$wb.Save()
